$wb = $excel.ActiveWorkbook

$wsJournal = $wb.Worksheets.Item("Journal")
$wsTotaux  = $wb.Worksheets.Item("Totaux")

# ---------------------------------------------------------------------------
# Journal sheet: append two new rows (22 & 23) to the "Tableau1" listobject,
# cloning the formatting of the last existing row (21) so no new cell
# styles are introduced.
# ---------------------------------------------------------------------------

$wsJournal.Range("A21:E21").Copy()
$wsJournal.Range("A22:E23").PasteSpecial(-4122)

$wsJournal.Range("A22").Value = 44980
$wsJournal.Range("B22").Value = 3
$wsJournal.Range("C22").Value = 0.03125
$wsJournal.Range("D22").Value = "Analyse"
$wsJournal.Range("E22").Value = "Analyse du tableau de bord du parc informatique "

$wsJournal.Range("A23").Value = 44980
$wsJournal.Range("B23").Value = 3
$wsJournal.Range("C23").Value = 0.03125
$wsJournal.Range("D23").Value = "Contrat"
$wsJournal.Range("E23").Value = "Rédaction du contrat d'altérnance pour le pré TPI et le TPI"

# Grow the Journal table ("Tableau1") to include the two new rows.
$loJournal = $wsJournal.ListObjects.Item("Tableau1")
$loJournal.Resize($wsJournal.Range("A1:E23"))

# ---------------------------------------------------------------------------
# Totaux sheet: insert a new weekly-total row before the grand "Total" row,
# pushing the Total row from row 9 down to row 10, and refresh the sums.
# ---------------------------------------------------------------------------

$wsTotaux.Rows(9).Insert()

$wsTotaux.Range("A8:B8").Copy()
$wsTotaux.Range("A9:B9").PasteSpecial(-4122)

$wsTotaux.Range("A9").Value = 44980
$wsTotaux.Range("B9").Formula = "=SUM(Journal!C22:C23)"

$wsTotaux.Range("B10").Formula = "=SUM(B2:B9)"

# Grow the Totaux table ("Tableau2") to include the new row.
$loTotaux = $wsTotaux.ListObjects.Item("Tableau2")
$loTotaux.Resize($wsTotaux.Range("A1:B10"))

# ---------------------------------------------------------------------------
# View state: Journal becomes the active sheet/tab, with a new selection on
# each sheet.
# ---------------------------------------------------------------------------

$wsTotaux.Range("F13").Select()

$wsJournal.Activate()
$wsJournal.Range("G25").Select()
